$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the Price column (D) as text so numeric-looking strings
# (e.g. "1.001", "40.81") are not silently converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.785.46"
$ws.Range("E2").Value = "  -2.58%  "
$ws.Range("D3").Value = "1.785.41"
$ws.Range("E3").Value = "  -2.21%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "309.73"
$ws.Range("E5").Value = "  -2.26%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("D7").Value = "0.5107"
$ws.Range("E7").Value = "  -1.74%  "
$ws.Range("D8").Value = "0.3874"
$ws.Range("E8").Value = "  -0.33%  "
$ws.Range("D9").Value = "0.07836"
$ws.Range("E9").Value = "  -6.64%  "
$ws.Range("D10").Value = "1.088"
$ws.Range("E10").Value = "  -2.46%  "
$ws.Range("D11").Value = "40.81"
$ws.Range("E11").Value = "  -2.66%  "
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("D13").Value = "6.214"
$ws.Range("E13").Value = "  -3.34%  "
$ws.Range("D14").Value = "20.17"
$ws.Range("E14").Value = "  -4.45%  "
$ws.Range("D15").Value = "1.779.06"
$ws.Range("E15").Value = "  -2.23%  "
$ws.Range("D16").Value = "7.211"
$ws.Range("E16").Value = "  -4.09%  "
$ws.Range("D17").Value = "91.27"
$ws.Range("E17").Value = "  -2.03%  "
$ws.Range("E18").Value = "  -4.85%  "
$ws.Range("D19").Value = "0.06520"
$ws.Range("E19").Value = "  -1.05%  "
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("D21").Value = "17.00"
$ws.Range("E21").Value = "  -4.46%  "
$ws.Range("D22").Value = "5.901"
$ws.Range("E22").Value = "  -2.80%  "
$ws.Range("D23").Value = "27.850.31"
$ws.Range("E23").Value = "  -2.46%  "
$ws.Range("D24").Value = "11.00"
$ws.Range("E24").Value = "  -3.76%  "
$ws.Range("D25").Value = "2.224"
$ws.Range("E25").Value = "  -2.61%  "
$ws.Range("D26").Value = "160.67"
$ws.Range("E26").Value = "  +0.91%  "
$ws.Range("D27").Value = "20.18"
$ws.Range("E27").Value = "  -4.30%  "
$ws.Range("D28").Value = "1.982.90"
$ws.Range("E28").Value = "  -2.41%  "
$ws.Range("D29").Value = "2.359"
$ws.Range("E29").Value = "  -1.75%  "
$ws.Range("D30").Value = "123.44"
$ws.Range("E30").Value = "  -1.78%  "
$ws.Range("D31").Value = "0.1077"
$ws.Range("E31").Value = "  -1.47%  "
$ws.Range("D32").Value = "1.032"
$ws.Range("E32").Value = "  -6.22%  "
$ws.Range("D33").Value = "3.632"
$ws.Range("E33").Value = "  -0.88%  "
$ws.Range("D34").Value = "5.478"
$ws.Range("E34").Value = "  -4.36%  "
$ws.Range("D35").Value = "0.07047"
$ws.Range("E35").Value = "  -5.19%  "
$ws.Range("D36").Value = "0.02302"
$ws.Range("E36").Value = "  -2.71%  "
$ws.Range("D37").Value = "8.772"
$ws.Range("E37").Value = "  -0.36%  "
$ws.Range("D38").Value = "0.2121"
$ws.Range("E38").Value = "  -4.52%  "
$ws.Range("D39").Value = "11.47"
$ws.Range("E39").Value = "  +0.39%  "
$ws.Range("D40").Value = "4.975"
$ws.Range("E40").Value = "  -4.81%  "
$ws.Range("D41").Value = "0.6077"
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("D43").Value = "1.148"
$ws.Range("E43").Value = "  -3.57%  "
$ws.Range("D44").Value = "1.314"
$ws.Range("E44").Value = "  -6.09%  "
$ws.Range("D45").Value = "13.07"
$ws.Range("E45").Value = "  -4.08%  "
$ws.Range("D46").Value = "0.5875"
$ws.Range("E46").Value = "  -1.51%  "
$ws.Range("E47").Value = "  -2.50%  "
$ws.Range("D48").Value = "124.37"
$ws.Range("E48").Value = "  -1.81%  "
$ws.Range("E49").Value = "  -0.47%  "
$ws.Range("D50").Value = "1.903"
$ws.Range("E50").Value = "  -4.46%  "
$ws.Range("D51").Value = "0.06819"
$ws.Range("E51").Value = "  -2.42%  "

# Restore default (General) formatting on column D so the saved file
# does not carry a residual number-format style on these cells.
$ws.Range("D2:D51").ClearFormats()
